# "alterações referente a caixa e financeiro"
#
# 1. Column F (col 6) becomes wider (to fit the new row-2 header text).
# 2. Row 2 gets taller to accommodate the now-wrapped header text.
# 3. The "Nº" column (C) is renumbered: row 3 -> 1, row 4 -> 2.
# 4. Cells E4/F4/J4 pick up the same background fill that already shows
#    on the corresponding cells in row 3 (E3/F3 -> light blue, J3 -> light
#    peach) - i.e. "paint" row 3's fill onto row 4's matching cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column F.
$ws.Columns.Item(6).ColumnWidth = 14.6

# Make row 2 taller.
$ws.Rows.Item(2).RowHeight = 26.25

# Renumber the "Nº" column.
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 2

# Apply the same fill as row 3 to the matching row-4 cells.
$ws.Range("E4:F4").Interior.Color = $ws.Range("E3").Interior.Color
$ws.Range("J4").Interior.Color = $ws.Range("J3").Interior.Color
